# "Make clear we don't know the distance"
#
# On the single slide of the deck, the text box named "ZoneTexte 32"
# (shape id 33) reads "d = 300cm". The known numeric value "300" is
# replaced with a placeholder "???" so the label becomes "d = ???cm",
# signalling that the distance is unknown.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item("ZoneTexte 32")
$tr = $shp.TextFrame.TextRange
$full = $tr.Text

$needle = "300"
$start = $full.IndexOf($needle)

if ($start -ge 0) {
    $sub = $tr.Characters($start + 1, $needle.Length)
    $sub.Text = "???"
}
